# SprintBacklog and Team member report final updates
#
# Appends "Created configuration management plan." to three team members'
# "role duties performed this week" cells (B7 on Yong, Vasilis and Brody),
# refreshes row heights for the now-longer wrapped text, and updates the
# saved selection/active-sheet state to match the author's final view
# (Brody's tab active, cursor on B16).

$wb = $excel.ActiveWorkbook

$suffix = " Created configuration management plan."

# --- Yong: row 7 (role duties performed this week) ---
$wsYong = $wb.Worksheets.Item("Yong")
$wsYong.Range("B7").Value = "Created GRL and UCM models document, Implemented front-end for route display." + $suffix
$wsYong.Rows.Item(7).RowHeight = 39
[void]$wsYong.Range("B7").Select()

# --- Vasilis: row 7 (role duties performed this week) ---
$wsVasilis = $wb.Worksheets.Item("Vasilis")
$wsVasilis.Range("B7").Value = "Created user interface for creating a schedule, adding classes, saving a schedule, recalling a schedule, schedule selection, and detail display. Updated SRS and UML document. Created GRM and UCM models document." + $suffix
$wsVasilis.Rows.Item(7).RowHeight = 78
[void]$wsVasilis.Range("B7").Select()

# --- Brody: row 7 (role duties performed this week) ---
$wsBrody = $wb.Worksheets.Item("Brody")
$wsBrody.Range("B7").Value = "Updated daily scrum report, Updated sprint backlog at beginning & end of sprint 3, Updates SRS and UML document, Completed sprint 3 review document, Updated team member report." + $suffix
$wsBrody.Rows.Item(7).RowHeight = 78

# --- Sakshyam: just reset the saved view state (loses tabSelected / topLeftCell scroll) ---
$wsSakshyam = $wb.Worksheets.Item("Sakshyam")
[void]$wsSakshyam.Range("B10").Select()

# --- Brody becomes the final active sheet/tab, with the cursor parked on B16 ---
[void]$wsBrody.Activate()
[void]$wsBrody.Range("B16").Select()
